# Automatische test-sync: 2025-07-29 21:57:50
# Appends the result of a new test mail (Testmail #14: CE-certificaten) as a
# new row at the bottom of the historical-responses log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A27:J27")

$newRow.Cells.Item(1, 1).Value  = "Testmail #14: Heb je de CE-certificaten van dit product?"
$newRow.Cells.Item(1, 2).Value  = "Beste klant,`nDank u voor uw e-mail. Wij kunnen u bevestigen dat dit product over de vereiste CE-certificaten beschikt. Mocht u nog verdere vragen hebben of meer informatie nodig hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$newRow.Cells.Item(1, 3).Value  = "Heb je de CE-certificaten van dit product?"
$newRow.Cells.Item(1, 4).Value  = "mailmind.test@zohomail.eu"
$newRow.Cells.Item(1, 5).Value  = "Productinformatie"
$newRow.Cells.Item(1, 6).Value  = "2025-07-29 21:57:40"
$newRow.Cells.Item(1, 7).Value  = "Ja"
$newRow.Cells.Item(1, 8).Value  = "Nee"
$newRow.Cells.Item(1, 9).Value  = "Ja"
$newRow.Cells.Item(1, 10).Value = "Nee"

# Writing multi-line content auto-expands the row height; reset it so the
# row matches the default (no explicit height) like every other data row.
$ws.Rows.Item(27).AutoFit()
